# Append the new "14-10-2025" gold-price row (row 23) to Sheet1,
# mirroring the existing Date/Gold-data rows (A = date, B = price text).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A23").Value = "14-10-2025"
$ws.Range("B23").Value = "The price of gold in India today is ₹12,868 per gram for 24 karat gold, ₹11,795 per gram for 22 karat gold and ₹9,651 per gram for 18 karat gold (also called 999 gold)."
